$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (date style) of the last existing row's date cell
# down into the new row so the new date cell keeps the same number format.
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(8, 1).Value = 42611.885752314818
$ws.Cells.Item(8, 2).Value = 0
$ws.Cells.Item(8, 3).Value = 54
$ws.Cells.Item(8, 4).Value = 43
$ws.Cells.Item(8, 5).Value = 45
$ws.Cells.Item(8, 6).Value = 54
$ws.Cells.Item(8, 7).Value = 13501
$ws.Cells.Item(8, 8).Value = 12545
$ws.Cells.Item(8, 9).Value = 2480
$ws.Cells.Item(8, 10).Value = 239
$ws.Cells.Item(8, 11).Value = 190
$ws.Cells.Item(8, 12).Value = 10
$ws.Cells.Item(8, 13).Value = 12
$ws.Cells.Item(8, 14).Value = "Noun"
